$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11, column B ("R40") is updated to the text value "1"
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = "1"
